$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "when was the last order i placed"
$ws.Range("A12").Value = "hello"
$ws.Range("A13").Value = "when was the last order i placed"
$ws.Range("A14").Value = "when was the last order i placed"
